$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author removed two rows (old rows 13 and 14, which held the two
# "Docentes responsaveis" name values with no label in column A) by
# deleting them outright. Everything below shifts up by two rows,
# carrying its row height along with it.
$ws.Rows("13:14").Delete()

# After the shift, a handful of the long descriptive paragraphs that used
# to sit next to several labels were cleared out and replaced with
# (duplicated) short values taken from elsewhere on the sheet, exactly as
# captured by the target workbook.
$ws.Range("B10").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("C10").Value = "3577649 - Carlos Angelo Nunes"

$ws.Range("B13").Value = "01/01/2022"
$ws.Range("C13").Value = "01/01/2022"

$ws.Range("B15").Value = "3577649 - Carlos Angelo Nunes"
$ws.Range("C15").Value = "3577649 - Carlos Angelo Nunes"

$ws.Range("B18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C18").Value = "519033 - Carlos Yujiro Shigue"

$ws.Range("B19").Value = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."
$ws.Range("C19").Value = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."

$ws.Range("B20").Value = "MF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo docente supervisor do estágio."
$ws.Range("C20").Value = "MF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo docente supervisor do estágio."

$ws.Range("B21").Value = "Não será oferecida recuperação."
$ws.Range("C21").Value = "Não será oferecida recuperação."
